# Add 2022-Q4 data
# -----------------------------------------------------------------
# 1) Update the "总计" (summary) sheet: insert a new row 2 holding the
#    2022-Q4 summary figures, pushing the existing quarters down.
# -----------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)

$summary.Rows(2).Insert()
$summary.Range("B2:D2").ClearFormats()

# Column A on the data rows carries the bordered/bold "index" style;
# grab it from row 3 (still intact) and stamp it onto the new row's A cell.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 13
$summary.Range("D2").Value = 7.84

# The "#" column is a 0-based running index; every pre-existing row shifts
# down one position so its index must be bumped by 1 as well.
for ($i = 3; $i -le 8; $i++) {
    $summary.Range("A" + $i).Value = $i - 2
}

# -----------------------------------------------------------------
# 2) Insert a brand-new "2022-Q4" sheet right after "总计", holding the
#    per-fund breakdown, matching the layout of the other quarter sheets.
# -----------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $summary)
$q4.Name = "2022-Q4"

$refSheet = $wb.Worksheets.Item("2022-Q3")
$refSheet.Range("A1:H2").Copy()
$q4.Range("A1:H2").PasteSpecial(-4122)

$q4.Range("A2:H2").Copy()
$q4.Range("A3:H14").PasteSpecial(-4122)

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$rows = @(
  @(0,  "206009", "鹏华新兴产业混合",             "43.49", "92.64", "7.81", "3.3966", 1),
  @(1,  "012930", "中庚价值先锋股票",             "68.71", "94.78", "4.23", "2.9064", 9),
  @(2,  "920003", "中金新锐股票A",                "19.30", "91.44", "3.47", "0.6697", 10),
  @(3,  "398001", "中海优质成长混合",             "11.10", "89.86", "3.03", "0.3363", 9),
  @(4,  "011542", "鹏华远见回报三年持有期混合",   "14.77", "60.61", "1.86", "0.2747", 10),
  @(5,  "920923", "中金新锐股票C",                "3.33",  "91.44", "3.47", "0.1156", 10),
  @(6,  "920002", "中金精选股票A",                "3.21",  "93.08", "2.77", "0.0889", 9),
  @(7,  "015133", "华安鼎安优选一年持有混合A",    "1.31",  "36.35", "1.02", "0.0134", 5),
  @(8,  "009719", "招商增浩一年定期开放混合C",    "0.71",  "21.56", "1.59", "0.0113", 2),
  @(9,  "161727", "招商增荣灵活配置混合（LOF）",  "0.52",  "55.51", "2.16", "0.0112", 1),
  @(10, "015134", "华安鼎安优选一年持有混合C",    "0.64",  "36.35", "1.02", "0.0065", 5),
  @(11, "009718", "招商增浩一年定期开放混合A",    "0.38",  "21.56", "1.59", "0.0060", 2),
  @(12, "920922", "中金精选股票C",                "0.12",  "93.08", "2.77", "0.0033", 9)
)

$r = 2
foreach ($row in $rows) {
    $textRange = "B" + $r + ":G" + $r
    $q4.Range($textRange).NumberFormat = "@"

    $q4.Range("A" + $r).Value = $row[0]
    $q4.Range("B" + $r).Value = $row[1]
    $q4.Range("C" + $r).Value = $row[2]
    $q4.Range("D" + $r).Value = $row[3]
    $q4.Range("E" + $r).Value = $row[4]
    $q4.Range("F" + $r).Value = $row[5]
    $q4.Range("G" + $r).Value = $row[6]
    $q4.Range("H" + $r).Value = $row[7]
    $r = $r + 1
}

# Strip the leftover "@" text-number-format from B:G (forced above only so
# the fund codes/figures are stored as text, matching the source data) by
# re-pasting the clean/default format from an always-blank helper cell.
for ($i = 2; $i -le 14; $i++) {
    $src = "I" + $i
    $dst = "B" + $i + ":G" + $i
    $q4.Range($src).Copy()
    $q4.Range($dst).PasteSpecial(-4122)
}
